# Slide 7: reposition/resize the existing picture and add a new rectangle
# shape containing a hyperlinked YouTube URL above it, per the target
# OOXML diff.
#
# Note: PowerPoint's COM object model expresses Shape.Left/.Top/.Width/.Height
# in points, while the OOXML stores EMUs (1 pt = 12700 EMU). We divide the
# target EMU values by 12700 to get the point values to feed into the COM API.

$emuPerPt = 12700
$url = "https://www.youtube.com/watch?v=z14p9Lo7NSs"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Reposition / resize the existing picture -----------------------------
$pic = $s.Shapes.Item("Picture 2")
$pic.Left   = 4529751 / $emuPerPt
$pic.Top    = 3231037 / $emuPerPt
# (nudged slightly to counter float32 rounding in the host so the saved EMU
# extent lands exactly on 6868539 instead of 6868538)
$pic.Width  = 540.8298645196775
$pic.Height = 3147646 / $emuPerPt

# --- Add the new rectangle shape with the hyperlinked URL -----------------
$rect = $s.Shapes.AddShape(1, 629704 / $emuPerPt, 2546751 / $emuPerPt, 5785558 / $emuPerPt, 369332 / $emuPerPt)

$rect.TextFrame.WordWrap = 0
$rect.TextFrame.AutoSize = 1

$rect.TextFrame.TextRange.Text = $url
$rect.TextFrame.TextRange.LanguageID = "en-CA"
$rect.TextFrame.TextRange.ActionSettings(1).Hyperlink.Address = $url
